# Sara Alert "Invalid Monitorees" template — add three new Race columns
# (Race Unknown, Race Other, Race Refused to Answer) to the header row,
# mirroring the formatting of the preceding header cell, then leave the
# selection on the newly added last header cell as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1 currently ends at CU1 ("Sexual Orientation"). Copy that
# cell's style across the three new header cells, then set their text.
$ws.Range("CU1").Copy($ws.Range("CV1:CX1"))

$ws.Range("CV1").Value = "Race Unknown"
$ws.Range("CW1").Value = "Race Other"
$ws.Range("CX1").Value = "Race Refused to Answer"

# Match the saved selection state left by the edit.
$ws.Range("CX6").Select() | Out-Null
